# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing N/O/P columns
# (Late / Date / Outstanding) one column to the right, and make the
# "Repayment schedule" sheet the active tab (it was previously
# "Transactions").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the
# freshly inserted column can inherit the same width, mirroring what
# Excel itself does when a column is inserted via the UI.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N; existing N,O,P shift right to O,P,Q.
$ws.Columns("N").Insert()

# Match the new column's width to its left neighbour (M).
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab, with the new
# selection positioned at R7 (one column to the right of the previous
# selection on this sheet, consistent with the inserted column).
$ws.Activate()
$ws.Range("R7").Select()
